# Product backlog update: logic wired up for the "stay" button.
#  - header "priority" column repurposed as "status"
#  - story points filled in for several backlog items
#  - BJ-003 (draw a card / get a score) marked as "finished"
#  - backlog items about the "dealer's" cards/score reworded to the
#    "other player's" cards/status (multiplayer support)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename the "priority" column to "status"
$ws.Range("F2").Value = "status"

# Story points (column E)
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 2
$ws.Range("E10").Value = 3

# BJ-003 is finished
$ws.Range("F5").Value = "finished"

# BJ-007: looking at the dealer's cards -> the other player's cards
$ws.Range("C9").Value = "be able to look at the other player's cards"

# BJ-010: seeing the dealer's score -> the other player's status
$ws.Range("C12").Value = "see the other player's status (score etc. as same as BJ-009)"

# Leave the selection where the user last clicked while making this edit
$ws.Range("E5").Select()
